# "Temps de jeu" workbook update: add a 5th CDF cup-match ("CDF T5") column
# group (minutes / T-R-NR-HG status / buts / passes D) and fill in the
# per-player data for the new match. Excel's own formulas (totals, COUNTIFs,
# etc.) recompute automatically once the raw cells are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: rename the new match-group header cell (was still showing the
# generic "CDF" shared string) to the distinct "CDF T5" label.
$ws.Range("HG1").Value = "CDF T5"

# Row 2 - Alban Rambaud: did not play (hors groupe -> NR this time)
$ws.Range("HH2").Value = "NR"

# Row 3 - Jassim Assoul: out of the squad (HG)
$ws.Range("HH3").Value = "HG"

# Row 4 - Enzo Vita: started, played full 90'
$ws.Range("HG4").Value = 90
$ws.Range("HH4").Value = "T"

# Row 5 - Romain Thunet: started, played 85'
$ws.Range("HG5").Value = 85
$ws.Range("HH5").Value = "T"

# Row 6 - Amine Taiar: out of the squad (HG)
$ws.Range("HH6").Value = "HG"

# Row 7 - Naim Ighbane: started, played full 90'
$ws.Range("HG7").Value = 90
$ws.Range("HH7").Value = "T"

# Row 8 - Hedi Nasri: out of the squad (HG)
$ws.Range("HH8").Value = "HG"

# Row 9 - Mattheo Haon: started, played full 90'
$ws.Range("HG9").Value = 90
$ws.Range("HH9").Value = "T"

# Row 10 - Maé Clavel: out of the squad (HG)
$ws.Range("HH10").Value = "HG"

# Row 11 - Levy Ndoutoume: out of the squad (HG)
$ws.Range("HH11").Value = "HG"

# Row 13 - Rayane Chayebi: out of the squad (HG)
$ws.Range("HH13").Value = "HG"

# Row 14 - Ilan Ihaddadene: came on as substitute, played 45', scored 1
$ws.Range("HG14").Value = 45
$ws.Range("HH14").Value = "R"
$ws.Range("HI14").Value = 1

# Row 15 - Karahali Souaré: started, played 45'
$ws.Range("HG15").Value = 45
$ws.Range("HH15").Value = "T"

# Row 16 - Amir Etien: came on as substitute, played 45'
$ws.Range("HG16").Value = 45
$ws.Range("HH16").Value = "R"

# Row 17 - Karim Belmahi: out of the squad (HG)
$ws.Range("HH17").Value = "HG"

# Row 18 - Emmanuel Valey: started, played full 90'
$ws.Range("HG18").Value = 90
$ws.Range("HH18").Value = "T"

# Row 19 - Jeremie Laurent: started, played 65'
$ws.Range("HG19").Value = 65
$ws.Range("HH19").Value = "T"

# Row 20 - Sofiane Belle: started, played 45'
$ws.Range("HG20").Value = 45
$ws.Range("HH20").Value = "T"

# Row 21 - Amir Kherrab: out of the squad (HG)
$ws.Range("HH21").Value = "HG"

# Row 22 - Naim Dhib: started, played full 90'
$ws.Range("HG22").Value = 90
$ws.Range("HH22").Value = "T"

# Row 23 - Wael Fareh: out of the squad (HG)
$ws.Range("HH23").Value = "HG"

# Row 24 - Yoan Zouma: came on as substitute, played 5'
$ws.Range("HG24").Value = 5
$ws.Range("HH24").Value = "R"

# Row 25 - Yoann Martelat: came on as substitute, played 25'
$ws.Range("HG25").Value = 25
$ws.Range("HH25").Value = "R"

# Row 26 - Omar Benyounes: out of the squad (HG)
$ws.Range("HH26").Value = "HG"

# Row 27 - Ilyes Boughanmi: started, played full 90'
$ws.Range("HG27").Value = 90
$ws.Range("HH27").Value = "T"

# Row 28 - Malik Boussaid: started, played full 90'
$ws.Range("HG28").Value = 90
$ws.Range("HH28").Value = "T"

# Row 29 - Kamal Bafounta: out of the squad (HG)
$ws.Range("HH29").Value = "HG"

# Leave the cursor where the author ended up editing.
$ws.Range("HK23").Select()
